$wb = $excel.ActiveWorkbook

# --- Sheet "Puntaje por Tareas" ---
$ws1 = $wb.Worksheets.Item("Puntaje por Tareas")

# Mark "Entrega en tiempo" (column D) as delivered on time (1) for the
# "Organizacion/Calendarizacion..." task (row 3) and the
# "Desarrollo de diseño/Mockups/WireFrames" task (row 5), matching the
# formatting already used by the other "delivered on time" cells (D4/D6).
$ws1.Range("D4").Copy()
$ws1.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("D3").Value = 1

$ws1.Range("D4").Copy()
$ws1.Range("D5").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("D5").Value = 1

# --- Sheet "Puntaje Promedio" ---
$ws2 = $wb.Worksheets.Item("Puntaje Promedio")

# Update the current points earned by each team member.
$ws2.Range("B3").Value = 2   # Rodrigo Euan
$ws2.Range("B4").Value = 2   # Javier Quijano
$ws2.Range("B7").Value = 2   # Gerardo Dueñas

$wb.Save()
